$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new product rows before the current row 9 ("جهاز محلول ") ---
$ws.Rows("9:10").Insert()

# Copy the formatting (styles / number formats / fonts / fills / borders) from
# the two existing product rows (7 and 8) onto the freshly inserted blank rows
# so the new rows look like every other product row in the table.
$ws.Range("A7:Q7").Copy()
$ws.Range("A9:Q9").PasteSpecial(-4122)
$ws.Range("A8:Q8").Copy()
$ws.Range("A10:Q10").PasteSpecial(-4122)

# Match row heights of the template rows
$ws.Rows("9").RowHeight = $ws.Rows("7").RowHeight
$ws.Rows("10").RowHeight = $ws.Rows("8").RowHeight

# Recreate the merged-cell layout used by every product row
$ws.Range("A9:B9").Merge()
$ws.Range("C9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("N9:O9").Merge()

$ws.Range("A10:B10").Merge()
$ws.Range("C10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:M10").Merge()
$ws.Range("N10:O10").Merge()

# --- Fill in the new product rows ---
# Row 9: ROWACHOL 45 CAPSULES
$ws.Range("A9").Value = 3
$ws.Range("C9").Value = "ROWACHOL 45 CAPSULES"
$ws.Range("H9").Value = "3:1"
$ws.Range("L9").Value = "1"
$ws.Range("N9").Value = "72.00"
$ws.Range("P9").Value = "23.7600"
$ws.Range("Q9").Value = "0:1"

# Row 10: بلاستر مترسيلك 2 سم
$ws.Range("A10").Value = 4
$ws.Range("C10").Value = "بلاستر مترسيلك 2 سم"
$ws.Range("H10").Value = "37:0"
$ws.Range("L10").Value = "0"
$ws.Range("N10").Value = "15.00"
$ws.Range("P10").Value = "15.0000"
$ws.Range("Q10").Value = "1:0"

# --- Renumber the rows that were pushed down (previously 9 & 10, now 11 & 12) ---
$ws.Range("A11").Value = 5
$ws.Range("A12").Value = 6

# --- Update the grand-total cell (previously P11, now P13) ---
$ws.Range("P13").Value = 163.49

# --- Update the generated-at timestamp (previously row 12, now row 14) ---
$ws.Range("A14").Value = "Friday, 29 August, 2025 11:49 AM"
